$d = $word.ActiveDocument

# =====================================================================
# Phase 1: text-content edits
# =====================================================================

# "value是-2.064" -> "value" + "的绝对值为" + "2.06" + "9"
$c = $d.Content
$c.Find.Execute("value是-2.064", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$isStart = $c.Start + 5
$isRun = $d.Range($isStart, $isStart + 1)
$isRun.Text = "的绝对值为"

$c2 = $d.Content
$c2.Find.Execute("-2.064", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$negRun = $d.Range($c2.Start, $c2.End)
$negRun.Text = "2.06"
$after206 = $c2.Start + 4

# grab formatting (no w:hint="eastAsia") from the still-intact "064" run via
# the single leading digit right after "大于2." and transplant it so the new
# "9" run gets the same (non-eastAsia-hinted) run properties.
$c3 = $d.Content
$c3.Find.Execute("大于2.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$noHintSrcStart = $c3.End
$noHintSrc = $d.Range($noHintSrcStart, $noHintSrcStart + 1)

$dst = $d.Range($after206, $after206)
$dst.FormattedText = $noHintSrc.FormattedText
$newDigit = $d.Range($after206, $after206 + 1)
$newDigit.Text = "9"

# second "2.064" -> "2.069" (only the "064" -> "069" part changes; reuse the
# existing un-hinted run in place so its formatting carries over untouched)
$c4 = $d.Content
$c4.Find.Execute("大于2.064", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$secondNum = $d.Range($c4.End - 3, $c4.End)
$secondNum.Text = "069"

# =====================================================================
# Phase 2: move the _GoBack bookmark from after "则拒绝原假设。" to between
# "值的绝对值" and "大于"
# =====================================================================

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$c5 = $d.Content
$c5.Find.Execute("值的绝对值大于", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mid = $c5.Start + 5
$bmr = $d.Range($mid, $mid)
$d.Bookmarks.Add("_GoBack", $bmr)

# =====================================================================
# Phase 3: re-establish the exact run boundaries the target XML expects.
# A text edit anywhere in the paragraph causes this engine to coalesce
# adjacent runs that happen to share identical formatting; toggling a
# (no-op) character property across a whole run's range forces the
# engine to keep/re-create a hard run boundary around it without
# altering the run's own formatting.
# =====================================================================

function Protect-Run([string]$needle, [int]$searchFrom) {
    $r = $d.Range($searchFrom, $d.Content.End)
    $r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rr = $d.Range($r.Start, $r.End)
    $rr.Font.Bold = $true
    $rr.Font.Bold = $false
    return $r.End
}

$p0 = $c.Start - 5   # a bit before "value", safely inside the same paragraph

$p = Protect-Run "value" $p0
$p = Protect-Run "的绝对值为" $p
$p = Protect-Run "2.06" $p
$p = Protect-Run "9" $p
$p = Protect-Run "，如果实验得到的" $p
$p = Protect-Run "t" $p
$p = Protect-Run "值的绝对值" $p
$p = Protect-Run "大于" $p
$p = Protect-Run "2." $p
$p = Protect-Run "069" $p
$p = Protect-Run "的概率小于" $p
$p = Protect-Run "0." $p
$p = Protect-Run "05" $p
$p = Protect-Run "则拒绝原假设。" $p

Write-Output $d.Content.Text
